# Edit script: update price/volume data (columns D-H) and normalize the
# "fixed_ticker" column (I) to "CLS" for all data rows (2-43) on Sheet1.
# This mirrors the commit "fixed issues restulted from extra files",
# which replaced stale per-row ticker strings in column I with the
# correct constant ticker ("CLS") and corrected the numeric OHLC/volume
# data in columns D:H.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Data table: Row, open_price(D), close_price(E), high_price(F), low_price(G), shares_outstanding(H)
$data = @(
    @(2, 11.06999969482422, 12.21000003814697, 12.44999980926514, 10.88000011444092, 115032686),
    @(3, 11.73999977111816, 13.38000011444092, 13.48999977111816, 11.13000011444092, 115032686),
    @(4, 12.94999980926514, 11.22000026702881, 13.35999965667725, 10.64999961853027, 115032686),
    @(5, 10.89000034332275, 9.060000419616699, 11.03999996185303, 8.069999694824219, 115032686),
    @(6, 10.86999988555908, 10.72999954223633, 11.35000038146973, 10.3100004196167, 115032686),
    @(7, 9.289999961853027, 11.07999992370606, 11.44999980926514, 8.829999923706055, 115032686),
    @(8, 10.85000038146973, 11.85000038146973, 12.19999980926514, 10.15999984741211, 115032686),
    @(9, 11.85999965667725, 13.88000011444092, 14.53999996185303, 11.60000038146973, 115032686),
    @(10, 14.52999973297119, 14.25, 14.73999977111816, 13.59000015258789, 115032686),
    @(11, 13.65999984741211, 11.89000034332275, 13.92000007629394, 11.85000038146973, 115032686),
    @(12, 12.40999984741211, 10.05000019073486, 12.76000022888184, 10.01000022888184, 115032686),
    @(13, 10.5, 10.10000038146973, 11.72000026702881, 10.05000019073486, 115032686),
    @(14, 10.31999969482422, 11.5, 11.61999988555908, 9.880000114440918, 115032686),
    @(15, 11.82999992370606, 11.8100004196167, 12.5, 11.68000030517578, 115032686),
    @(16, 10.93000030517578, 10.38000011444092, 11.21000003814697, 9.109999656677246, 115032686),
    @(17, 8.699999809265137, 9.930000305175779, 9.960000038146973, 8.630000114440918, 115032686),
    @(18, 8.520000457763672, 7.130000114440918, 8.970000267028809, 6.900000095367432, 115032686),
    @(19, 6.920000076293945, 7.090000152587891, 7.25, 6.409999847412109, 115032686),
    @(20, 7.159999847412109, 7.21999979019165, 7.28000020980835, 5.949999809265137, 115032686),
    @(21, 8.310000419616699, 9.039999961853027, 9.289999961853027, 8.090000152587891, 115032686),
    @(22, 3.349999904632568, 6.230000019073486, 6.519999980926514, 3.200000047683716, 115032686),
    @(23, 6.840000152587891, 8.270000457763672, 8.399999618530273, 6.400000095367432, 115032686),
    @(24, 6.889999866485596, 5.869999885559082, 7.71999979019165, 5.769999980926514, 115032686),
    @(25, 8.149999618530273, 8.069999694824219, 8.949999809265137, 7.820000171661377, 115032686),
    @(26, 8.390000343322754, 8.329999923706055, 8.909999847412109, 7.849999904632568, 115032686),
    @(27, 7.900000095367432, 8.869999885559082, 9.050000190734863, 7.010000228881836, 115032686),
    @(28, 8.939999580383301, 9.819999694824221, 10.52999973297119, 8.859999656677246, 115032686),
    @(29, 11.14000034332275, 12.44999980926514, 12.52999973297119, 10.14999961853027, 115032686),
    @(30, 11.94999980926514, 11.22999954223633, 12.02000045776367, 10.44999980926514, 115032686),
    @(31, 9.590000152587891, 10.53999996185303, 10.90999984741211, 9.119999885559082, 115032686),
    @(32, 8.529999732971191, 10.96000003814697, 11.68000030517578, 8.340000152587891, 115032686),
    @(33, 11.31999969482422, 13.32999992370606, 14.27999973297119, 10.98999977111816, 115032686),
    @(34, 12.84000015258789, 10.85999965667725, 13.01000022888184, 10.75, 115032686),
    @(35, 14.5, 21.97999954223633, 22.01000022888184, 14.36999988555908, 115032686),
    @(36, 24.5, 23.35000038146973, 27.57999992370605, 20.8700008392334, 115032686),
    @(37, 28.63999938964844, 34.36000061035156, 35.22999954223633, 26.6200008392334, 115032686),
    @(38, 45.09000015258789, 43.33000183105469, 51.11999893188477, 40.31999969482422, 115032686),
    @(39, 58, 52.43999862670898, 63.4900016784668, 48.34000015258789, 115032686),
    @(40, 50.93999862670898, 68.40000152587891, 70.80999755859375, 48.97000122070312, 115032686),
    @(41, 92.83000183105467, 123.4700012207031, 133.6699981689453, 86.12000274658203, 115032686),
    @(42, 79.02999877929688, 85.34999847412109, 95.22000122070312, 58.04999923706055, 115032686),
    @(43, 154.8999938964844, 199.8600006103516, 214.4700012207031, 144.2700042724609, 115032686)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 4).Value = $row[1]   # D - open_price
    $ws.Cells.Item($r, 5).Value = $row[2]   # E - close_price
    $ws.Cells.Item($r, 6).Value = $row[3]   # F - high_price
    $ws.Cells.Item($r, 7).Value = $row[4]   # G - low_price
    $ws.Cells.Item($r, 8).Value = $row[5]   # H - shares_outstanding
    $ws.Cells.Item($r, 9).Value = "CLS"     # I - fixed_ticker
}
